# Update the hidden "__footings__" worksheet: replace bracketed mapping-key
# labels in column C (e.g. "[key]", "[a]", "[('tuple', 'key')]") with
# slash-delimited equivalents (e.g. "/key/", "/a/", "/('tuple', 'key')/").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("__footings__")

$rowsKey   = 18..27
$rowsTuple = 28..29
$rowsA     = 43,44,47,48
$rowsB     = 45,46,49,50

foreach ($r in $rowsKey) {
    $ws.Cells.Item($r, 3).Value = "/key/"
}

foreach ($r in $rowsTuple) {
    $ws.Cells.Item($r, 3).Value = "/('tuple', 'key')/"
}

foreach ($r in $rowsA) {
    $ws.Cells.Item($r, 3).Value = "/a/"
}

foreach ($r in $rowsB) {
    $ws.Cells.Item($r, 3).Value = "/b/"
}
